$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: actualizar mensaje de conversión del día ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$conversionText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.9 = 40584.26 pesos`n✅ 40584.26 pesos = 9.84 = 934.77 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $conversionText

# --- tasas: actualizar tasas automáticas ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 101
$wsTasas.Range("O10").Value = 4099.01
$wsTasas.Range("N12").Value = 4125
$wsTasas.Range("O12").Value = 95.01000000000001
